$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new data row appended after existing row 5
$ws.Cells.Item(6, 1).Value = 42604.890810185185
$ws.Cells.Item(6, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(6, 2).Value = "Bag"
$ws.Cells.Item(6, 3).Value = 5007
$ws.Cells.Item(6, 4).Value = 6181
$ws.Cells.Item(6, 5).Value = 824
$ws.Cells.Item(6, 6).Value = 91
$ws.Cells.Item(6, 7).Value = 51
$ws.Cells.Item(6, 8).Value = 63
$ws.Cells.Item(6, 9).Value = 35
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(6, 11).Value = 9
$ws.Cells.Item(6, 12).Value = 18
$ws.Cells.Item(6, 13).Value = 81
